# Port battle calculator correction
# - Reorders the "Shallow water port" ship list (rows 4-21, column B) to the
#   corrected sequence (this also reshuffles the shared-string table the same
#   way the authoritative diff does).
# - Updates the BR values (column C) that changed for several ships.
# - Fixes the BR-total SUM ranges on both sheets so they include the last
#   data row.

$wb = $excel.ActiveWorkbook

$wsDeep = $wb.Worksheets.Item("Deep water port")
$wsShallow = $wb.Worksheets.Item("Shallow water port")

# --- Deep water port: extend SUM ranges to include row 35 ---
$wsDeep.Range("D3").Formula = "=SUM(D4:D35)"
$wsDeep.Range("E3").Formula = "=SUM(E4:E35)"

# --- Shallow water port: extend SUM ranges to include row 21 ---
$wsShallow.Range("D3").Formula = "=SUM(D4:D21)"
$wsShallow.Range("E3").Formula = "=SUM(E4:E21)"

# --- Shallow water port: corrected ship order (rows 4-21) ---
$shipNames = @(
    "Hercules",
    "Pandora",
    "Mercury",
    "Mortar Brig",
    "NavyBrig",
    "Niagara",
    "Prince de Neufchatel",
    "Rattlesnake",
    "Rattlesnake Heavy",
    "Snow",
    "Brig",
    "Pickle",
    "Cutter",
    "GunBoat",
    "Lynx",
    "Privateer",
    "Yacht",
    "Yacht Silver"
)

for ($i = 0; $i -lt $shipNames.Length; $i++) {
    $row = 4 + $i
    $wsShallow.Range("B$row").Value = $shipNames[$i]
}

# --- Shallow water port: corrected BR values (column C) ---
$wsShallow.Range("C4").Value = 100
$wsShallow.Range("C5").Value = 100
$wsShallow.Range("C13").Value = 80
$wsShallow.Range("C14").Value = 70
$wsShallow.Range("C15").Value = 55
$wsShallow.Range("C21").Value = 50
